$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 27289
$ws.Range("B2").Value = 43793
$ws.Range("C2").Value = 71082
$ws.Range("D2").Value = 1256.962779281228
$ws.Range("E2").Value = -617.2640911946989
$ws.Range("F2").Value = 0.1280079258180017
$ws.Range("G2").Value = -0.06128177334911047
$ws.Range("H2").Value = 0.01138819373920035
$ws.Range("I2").Value = 34301257.2838055
$ws.Range("J2").Value = -27031846.34568904
$ws.Range("K2").Value = 0.3839087251343519
$ws.Range("L2").Value = 2.036345216272679
$ws.Range("M2").Value = 1.268920252251869
$ws.Range("N2").Value = 7269410.938116461
